$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the "Price" / "Volume(1h)" columns with the latest scrape.
# Rows 14/15 also swap rank: TRON now outranks Litecoin.
#
# The Price column stores plain text (e.g. "1.001", "26.696.08") rather than
# numbers, so force text formatting before writing to stop Excel's
# auto-coercion from turning numeric-looking strings into real numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '26.696.08'
$ws.Range("E2").Value = '  -1.22%  '
$ws.Range("D3").Value = '1.796.33'
$ws.Range("E3").Value = '  -1.20%  '
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = '309.20'
$ws.Range("E5").Value = '  -0.42%  '
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  -0.01%  '
$ws.Range("D7").Value = '0.4433'
$ws.Range("E7").Value = '  +5.27%  '
$ws.Range("D8").Value = '0.3681'
$ws.Range("E8").Value = '  +0.57%  '
$ws.Range("D9").Value = '0.07354'
$ws.Range("E9").Value = '  +2.14%  '
$ws.Range("D10").Value = '0.8603'
$ws.Range("E10").Value = '  +2.24%  '
$ws.Range("D11").Value = '20.65'
$ws.Range("E11").Value = '  -0.72%  '
$ws.Range("D12").Value = '1.798.55'
$ws.Range("E12").Value = '  -0.99%  '
$ws.Range("D13").Value = '6.621'
$ws.Range("E13").Value = '  +0.03%  '
$ws.Range("B14").Value = 'TRON'
$ws.Range("C14").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D14").Value = '0.07066'
$ws.Range("E14").Value = '  +0.02%  '
$ws.Range("B15").Value = 'Litecoin'
$ws.Range("C15").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D15").Value = '91.85'
$ws.Range("E15").Value = '  +3.36%  '
$ws.Range("D16").Value = '5.270'
$ws.Range("E16").Value = '  +0.06%  '
$ws.Range("D17").Value = '1.001'
$ws.Range("E17").Value = '  -0.15%  '
$ws.Range("D18").Value = '0.000008686'
$ws.Range("E18").Value = '  -1.29%  '
$ws.Range("D19").Value = '1.001'
$ws.Range("E19").Value = '  -0.03%  '
$ws.Range("D20").Value = '14.81'
$ws.Range("E20").Value = '  -0.86%  '
$ws.Range("D21").Value = '26.716.58'
$ws.Range("E21").Value = '  -1.54%  '
$ws.Range("D22").Value = '5.162'
$ws.Range("E22").Value = '  +1.02%  '
$ws.Range("D23").Value = '10.81'
$ws.Range("E23").Value = '  +0.17%  '
$ws.Range("D24").Value = '1.978'
$ws.Range("E24").Value = '  +0.08%  '
$ws.Range("D25").Value = '152.03'
$ws.Range("E25").Value = '  +0.35%  '
$ws.Range("D26").Value = '2.180'
$ws.Range("E26").Value = '  -2.23%  '
$ws.Range("D27").Value = '18.43'
$ws.Range("E27").Value = '  +1.07%  '
$ws.Range("D28").Value = '5.182'
$ws.Range("E28").Value = '  -0.37%  '
$ws.Range("D29").Value = '117.18'
$ws.Range("D30").Value = '0.08776'
$ws.Range("E30").Value = '  -0.05%  '
$ws.Range("D31").Value = '0.7397'
$ws.Range("E31").Value = '  +0.22%  '
$ws.Range("D32").Value = '1.155'
$ws.Range("E32").Value = '  -1.61%  '
$ws.Range("D33").Value = '4.445'
$ws.Range("E33").Value = '  +0.83%  '
$ws.Range("D34").Value = '2.902'
$ws.Range("E34").Value = '  -1.72%  '
$ws.Range("E35").Value = '  -0.02%  '
$ws.Range("E36").Value = '  -0.78%  '
$ws.Range("D37").Value = '0.01956'
$ws.Range("E37").Value = '  +0.08%  '
$ws.Range("D38").Value = '0.05191'
$ws.Range("E38").Value = '  -0.68%  '
$ws.Range("D39").Value = '0.5251'
$ws.Range("E39").Value = '  +4.45%  '
$ws.Range("D40").Value = '2.826'
$ws.Range("E40").Value = '  -1.78%  '
$ws.Range("D41").Value = '6.969'
$ws.Range("E41").Value = '  -3.99%  '
$ws.Range("D42").Value = '0.1683'
$ws.Range("E42").Value = '  +0.00%  '
$ws.Range("D43").Value = '0.5056'
$ws.Range("E43").Value = '  +6.87%  '
$ws.Range("D44").Value = '8.449'
$ws.Range("E44").Value = '  -1.29%  '
$ws.Range("D45").Value = '1.977'
$ws.Range("E45").Value = '  +5.17%  '
$ws.Range("D46").Value = '10.44'
$ws.Range("E46").Value = '  -0.58%  '
$ws.Range("D47").Value = '105.06'
$ws.Range("E47").Value = '  -0.97%  '
$ws.Range("D48").Value = '1.000'
$ws.Range("E48").Value = '  -0.04%  '
$ws.Range("D49").Value = '1.665'
$ws.Range("E49").Value = '  +1.36%  '
$ws.Range("D50").Value = '0.06289'
$ws.Range("E50").Value = '  -1.10%  '
$ws.Range("D51").Value = '0.9162'
$ws.Range("E51").Value = '  +1.73%  '

# Restore default (unstyled) formatting on the D column so no stray number-format style remains.
$ws.Range("D2:D51").Style = "Normal"
